$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 18996
$ws.Range("J3").Value = 18996
$ws.Range("L3").Value = 18996
$ws.Range("N3").Value = -19224
$ws.Range("H32").Value = 994
$ws.Range("J32").Value = 992.75
$ws.Range("L32").Value = 992.75
$ws.Range("N32").Value = -1644.75
$ws.Range("H102").Value = 18996
$ws.Range("J102").Value = 18996
$ws.Range("L102").Value = 18996
$ws.Range("N102").Value = -25486
$ws.Range("H113").Value = 3663.25
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H131").Value = 5471.5
$ws.Range("I131").Value = 4443.5
$ws.Range("J131").Value = 6499.5
$ws.Range("K131").Value = 13330.5
$ws.Range("L131").Value = 19498.5
$ws.Range("M131").Value = -8290.5
$ws.Range("N131").Value = -29578.5
$ws.Range("H137").Value = 5568165
$ws.Range("I137").Value = 8334413.5
$ws.Range("K137").Value = 25003240.5
$ws.Range("M137").Value = -25000690.5
$ws.Range("H141").Value = 5065
$ws.Range("I141").Value = 2294.5
$ws.Range("J141").Value = 29999.5
$ws.Range("K141").Value = 6883.5
$ws.Range("L141").Value = 89998.5
$ws.Range("M141").Value = -1703.5
$ws.Range("N141").Value = -100358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1058.625
$ws.Range("I2").Value = 724.7778
$ws.Range("J2").Value = 1487.8572
$ws.Range("K2").Value = 724.7778
$ws.Range("L2").Value = 1487.8572
$ws.Range("M2").Value = -611.7778
$ws.Range("N2").Value = -1713.8572
$ws.Range("H32").Value = 4318.0225
$ws.Range("I32").Value = 3794.7632
$ws.Range("K32").Value = 3794.7632
$ws.Range("M32").Value = -3507.7632
$ws.Range("H45").Value = 106602
$ws.Range("I45").Value = 106602
$ws.Range("K45").Value = 106602
$ws.Range("M45").Value = -106225
$ws.Range("H61").Value = 4564.625
$ws.Range("I61").Value = 2638.182
$ws.Range("J61").Value = 8802.799999999999
$ws.Range("K61").Value = 2638.182
$ws.Range("L61").Value = 8802.799999999999
$ws.Range("M61").Value = -2426.182
$ws.Range("N61").Value = -9226.799999999999
$ws.Range("H74").Value = 207734.03
$ws.Range("J74").Value = 3397.5
$ws.Range("L74").Value = 3397.5
$ws.Range("N74").Value = -5145.5
$ws.Range("H77").Value = 207734.03
$ws.Range("J77").Value = 3397.5
$ws.Range("L77").Value = 16987.5
$ws.Range("N77").Value = -25723.5
$ws.Range("H80").Value = 84705
$ws.Range("J80").Value = 84705
$ws.Range("L80").Value = 84705
$ws.Range("N80").Value = -86701
$ws.Range("H83").Value = 84705
$ws.Range("J83").Value = 84705
$ws.Range("L83").Value = 254115
$ws.Range("N83").Value = -264099
$ws.Range("H116").Value = 1058.625
$ws.Range("I116").Value = 724.7778
$ws.Range("J116").Value = 1487.8572
$ws.Range("K116").Value = 724.7778
$ws.Range("L116").Value = 1487.8572
$ws.Range("M116").Value = 1569.2222
$ws.Range("N116").Value = -6075.8572
$ws.Range("H122").Value = 2593.1282
$ws.Range("I122").Value = 2389.7778
$ws.Range("K122").Value = 7169.3334
$ws.Range("M122").Value = -4719.3334
$ws.Range("H132").Value = 2673.0334
$ws.Range("I132").Value = 1702.7368
$ws.Range("J132").Value = 4349
$ws.Range("K132").Value = 5108.2104
$ws.Range("L132").Value = 13047
$ws.Range("M132").Value = -2578.2104
$ws.Range("N132").Value = -18107
$ws.Range("H136").Value = 4564.625
$ws.Range("I136").Value = 2638.182
$ws.Range("J136").Value = 8802.799999999999
$ws.Range("K136").Value = 7914.545999999999
$ws.Range("L136").Value = 26408.4
$ws.Range("M136").Value = -5364.545999999999
$ws.Range("N136").Value = -31508.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1058.625
$ws.Range("I3").Value = 724.7778
$ws.Range("J3").Value = 1487.8572
$ws.Range("K3").Value = 724.7778
$ws.Range("L3").Value = 1487.8572
$ws.Range("M3").Value = -610.7778
$ws.Range("N3").Value = -1715.8572
$ws.Range("H44").Value = 21500
$ws.Range("J44").Value = 21500
$ws.Range("L44").Value = 21500
$ws.Range("N44").Value = -22494
$ws.Range("H54").Value = 77500.5
$ws.Range("I54").Value = 77500.5
$ws.Range("K54").Value = 77500.5
$ws.Range("M54").Value = -77016.5
$ws.Range("H99").Value = 5485.4287
$ws.Range("I99").Value = 3999.75
$ws.Range("K99").Value = 3999.75
$ws.Range("M99").Value = -2501.75
$ws.Range("H107").Value = 1425.6875
$ws.Range("I107").Value = 1246.0869
$ws.Range("J107").Value = 1884.6666
$ws.Range("K107").Value = 1246.0869
$ws.Range("L107").Value = 1884.6666
$ws.Range("M107").Value = 673.9131
$ws.Range("N107").Value = -5724.6666
$ws.Range("H134").Value = 3691.1287
$ws.Range("I134").Value = 3462.544
$ws.Range("J134").Value = 4693.385
$ws.Range("K134").Value = 10387.632
$ws.Range("L134").Value = 14080.155
$ws.Range("M134").Value = -7852.632
$ws.Range("N134").Value = -19150.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 13502.857
$ws.Range("J41").Value = 24840
$ws.Range("L41").Value = 24840
$ws.Range("N41").Value = -25696
$ws.Range("H59").Value = 69541.836
$ws.Range("I59").Value = 57500
$ws.Range("J59").Value = 75562.75
$ws.Range("K59").Value = 57500
$ws.Range("L59").Value = 75562.75
$ws.Range("M59").Value = -56355
$ws.Range("N59").Value = -77852.75
$ws.Range("H109").Value = 49499.5
$ws.Range("J109").Value = 49499.5
$ws.Range("L109").Value = 49499.5
$ws.Range("N109").Value = -51579.5
$ws.Range("H132").Value = 2330.1282
$ws.Range("I132").Value = 1876.381
$ws.Range("K132").Value = 5629.143
$ws.Range("M132").Value = -3099.143
$ws.Range("H134").Value = 2228.7222
$ws.Range("I134").Value = 2156.1667
$ws.Range("J134").Value = 2591.5
$ws.Range("K134").Value = 6468.500100000001
$ws.Range("L134").Value = 7774.5
$ws.Range("M134").Value = -3933.500100000001
$ws.Range("N134").Value = -12844.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 275037.34
$ws.Range("I9").Value = 275037.34
$ws.Range("K9").Value = 825112.02
$ws.Range("M9").Value = -824888.02
$ws.Range("H25").Value = 6954.1816
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 7549.6
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 22648.8
$ws.Range("M25").Value = -2831
$ws.Range("N25").Value = -22986.8
$ws.Range("H30").Value = 6954.1816
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 7549.6
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 22648.8
$ws.Range("M30").Value = -2898
$ws.Range("N30").Value = -22852.8
$ws.Range("H131").Value = 9782.789000000001
$ws.Range("I131").Value = 22422.285
$ws.Range("J131").Value = 2409.75
$ws.Range("K131").Value = 67266.855
$ws.Range("L131").Value = 7229.25
$ws.Range("M131").Value = -62226.855
$ws.Range("N131").Value = -17309.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 665
$ws.Range("I2").Value = 942
$ws.Range("J2").Value = 111
$ws.Range("K2").Value = 942
$ws.Range("L2").Value = 111
$ws.Range("M2").Value = -829
$ws.Range("N2").Value = -337
$ws.Range("H102").Value = 1039.0454
$ws.Range("I102").Value = 695
$ws.Range("K102").Value = 695
$ws.Range("M102").Value = 927
$ws.Range("H123").Value = 57999.668
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -94899
$ws.Range("H132").Value = 6945
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6945
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 20835
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -25895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1881.25
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("N50").Value = -51274
$ws.Range("H55").Value = 329
$ws.Range("J55").Value = 414
$ws.Range("L55").Value = 414
$ws.Range("N55").Value = -760
$ws.Range("H61").Value = 2254.55
$ws.Range("I61").Value = 2127
$ws.Range("K61").Value = 2127
$ws.Range("M61").Value = -1925
$ws.Range("H82").Value = 416.81818
$ws.Range("I82").Value = 375.57144
$ws.Range("K82").Value = 375.57144
$ws.Range("M82").Value = -14.57144
$ws.Range("H85").Value = 416.81818
$ws.Range("I85").Value = 375.57144
$ws.Range("K85").Value = 375.57144
$ws.Range("M85").Value = 872.4285600000001
$ws.Range("H113").Value = 2254.55
$ws.Range("I113").Value = 2127
$ws.Range("K113").Value = 2127
$ws.Range("M113").Value = 43
$ws.Range("H122").Value = 5999
$ws.Range("I122").Value = 5999
$ws.Range("K122").Value = 17997
$ws.Range("M122").Value = -15547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11124.25
$ws.Range("I62").Value = 8250
$ws.Range("J62").Value = 13998.5
$ws.Range("K62").Value = 8250
$ws.Range("L62").Value = 13998.5
$ws.Range("M62").Value = -7626
$ws.Range("N62").Value = -15246.5
$ws.Range("H65").Value = 11124.25
$ws.Range("I65").Value = 8250
$ws.Range("J65").Value = 13998.5
$ws.Range("K65").Value = 41250
$ws.Range("L65").Value = 69992.5
$ws.Range("M65").Value = -38130
$ws.Range("N65").Value = -76232.5
